$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 294.03125
$ws.Range("I19").Value = 278.66666
$ws.Range("J19").Value = 307.58823
$ws.Range("K19").Value = 278.66666
$ws.Range("L19").Value = 307.58823
$ws.Range("M19").Value = -103.66666
$ws.Range("N19").Value = -657.5882300000001
# Row 113
$ws.Range("H113").Value = 1500
$ws.Range("I113").Value = 1500
$ws.Range("K113").Value = 1500
$ws.Range("M113").Value = 1754
# Row 116
$ws.Range("H116").Value = 3531.7144
$ws.Range("J116").Value = 3876.4443
$ws.Range("L116").Value = 3876.4443
$ws.Range("N116").Value = -10760.4443
# Row 132
$ws.Range("H132").Value = 6066825.5
$ws.Range("I132").Value = 7411448.5
$ws.Range("K132").Value = 22234345.5
$ws.Range("M132").Value = -22231815.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 4
$ws.Range("H4").Value = 675
$ws.Range("I4").Value = 400
$ws.Range("J4").Value = 881.25
$ws.Range("K4").Value = 400
$ws.Range("L4").Value = 881.25
$ws.Range("M4").Value = -284
$ws.Range("N4").Value = -1113.25
# Row 32
$ws.Range("H32").Value = 6904.9785
$ws.Range("I32").Value = 5623.275
$ws.Range("K32").Value = 5623.275
$ws.Range("M32").Value = -5336.275
# Row 63
$ws.Range("H63").Value = 34484804
$ws.Range("I63").Value = 1863.1765
$ws.Range("J63").Value = 83335630
$ws.Range("K63").Value = 1863.1765
$ws.Range("L63").Value = 83335630
$ws.Range("M63").Value = -1177.1765
$ws.Range("N63").Value = -83337002
# Row 66
$ws.Range("H66").Value = 34484804
$ws.Range("I66").Value = 1863.1765
$ws.Range("J66").Value = 83335630
$ws.Range("K66").Value = 9315.8825
$ws.Range("L66").Value = 416678150
$ws.Range("M66").Value = -5883.8825
$ws.Range("N66").Value = -416685014
# Row 74
$ws.Range("H74").Value = 2875.3157
$ws.Range("I74").Value = 2377.1538
$ws.Range("J74").Value = 3954.6667
$ws.Range("K74").Value = 2377.1538
$ws.Range("L74").Value = 3954.6667
$ws.Range("M74").Value = -1503.1538
$ws.Range("N74").Value = -5702.6667
# Row 77
$ws.Range("H77").Value = 2875.3157
$ws.Range("I77").Value = 2377.1538
$ws.Range("J77").Value = 3954.6667
$ws.Range("K77").Value = 11885.769
$ws.Range("L77").Value = 19773.3335
$ws.Range("M77").Value = -7517.769
$ws.Range("N77").Value = -28509.3335
# Row 92
$ws.Range("H92").Value = 1679666.6
$ws.Range("J92").Value = 1679666.6
$ws.Range("L92").Value = 1679666.6
$ws.Range("N92").Value = -1684658.6
# Row 102
$ws.Range("H102").Value = 10418087
$ws.Range("I102").Value = 13890091
$ws.Range("J102").Value = 2074.75
$ws.Range("K102").Value = 13890091
$ws.Range("L102").Value = 2074.75
$ws.Range("M102").Value = -13888469
$ws.Range("N102").Value = -5318.75
# Row 110
$ws.Range("H110").Value = 1594.75
$ws.Range("I110").Value = 1134.5555
$ws.Range("J110").Value = 2975.3333
$ws.Range("K110").Value = 1134.5555
$ws.Range("L110").Value = 2975.3333
$ws.Range("M110").Value = 910.4445000000001
$ws.Range("N110").Value = -7065.3333
# Row 122
$ws.Range("H122").Value = 1985.2
$ws.Range("I122").Value = 1679
$ws.Range("J122").Value = 2444.5
$ws.Range("K122").Value = 5037
$ws.Range("L122").Value = 7333.5
$ws.Range("M122").Value = -2587
$ws.Range("N122").Value = -12233.5
# Row 134
$ws.Range("H134").Value = 31936
$ws.Range("J134").Value = 31936
$ws.Range("L134").Value = 31936
$ws.Range("N134").Value = -42076

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 58824812
$ws.Range("I105").Value = 66667720
$ws.Range("K105").Value = 66667720
$ws.Range("M105").Value = -66665973

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 25
$ws.Range("H25").Value = 2066.6667
$ws.Range("I25").Value = 1100
$ws.Range("K25").Value = 1100
$ws.Range("M25").Value = -926
# Row 105
$ws.Range("H105").Value = 1116.125
$ws.Range("I105").Value = 1074
$ws.Range("K105").Value = 1074
$ws.Range("M105").Value = 673
# Row 122
$ws.Range("H122").Value = 730.05554
$ws.Range("I122").Value = 723.13336
$ws.Range("J122").Value = 764.6667
$ws.Range("K122").Value = 2169.40008
$ws.Range("L122").Value = 2294.0001
$ws.Range("M122").Value = 280.5999199999997
$ws.Range("N122").Value = -7194.0001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 33
$ws.Range("H33").Value = 398.66666
$ws.Range("I33").Value = 222.5
$ws.Range("K33").Value = 1335
$ws.Range("M33").Value = -1052
# Row 55
$ws.Range("H55").Value = 1727.8235
$ws.Range("J55").Value = 1727.8235
$ws.Range("L55").Value = 5183.470499999999
$ws.Range("N55").Value = -5537.470499999999
# Row 92
$ws.Range("H92").Value = 270.70834
$ws.Range("I92").Value = 246.88235
$ws.Range("J92").Value = 328.57144
$ws.Range("K92").Value = 740.64705
$ws.Range("L92").Value = 985.71432
$ws.Range("M92").Value = 507.35295
$ws.Range("N92").Value = -3481.71432
# Row 95
$ws.Range("H95").Value = 25199.8
$ws.Range("J95").Value = 25199.8
$ws.Range("L95").Value = 75599.39999999999
$ws.Range("N95").Value = -79717.39999999999
# Row 107
$ws.Range("H107").Value = 3695.5
$ws.Range("J107").Value = 4335.4
$ws.Range("L107").Value = 13006.2
$ws.Range("N107").Value = -16846.2
# Row 113
$ws.Range("H113").Value = 661.9737
$ws.Range("I113").Value = 585.4545000000001
$ws.Range("J113").Value = 693.14813
$ws.Range("K113").Value = 1756.3635
$ws.Range("L113").Value = 2079.44439
$ws.Range("M113").Value = 413.6364999999998
$ws.Range("N113").Value = -6419.444390000001
# Row 131
$ws.Range("H131").Value = 30306870
$ws.Range("I131").Value = 83333790
$ws.Range("J131").Value = 5770.5713
$ws.Range("K131").Value = 250001370
$ws.Range("L131").Value = 17311.7139
$ws.Range("M131").Value = -249996330
$ws.Range("N131").Value = -27391.7139

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 57
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()
# Row 70
$ws.Range("H70").Value = 40912388
$ws.Range("J70").Value = 40003080
$ws.Range("L70").Value = 40003080
$ws.Range("N70").Value = -40003620
# Row 73
$ws.Range("H73").Value = 40912388
$ws.Range("J73").Value = 40003080
$ws.Range("L73").Value = 40003080
$ws.Range("N73").Value = -40004952
# Row 97
$ws.Range("H97").Value = 549.375
$ws.Range("I97").Value = 542.5
$ws.Range("J97").Value = 570
$ws.Range("K97").Value = 542.5
$ws.Range("L97").Value = 570
$ws.Range("M97").Value = -46.5
$ws.Range("N97").Value = -1562
# Row 113
$ws.Range("H113").Value = 1645.0454
$ws.Range("I113").Value = 1604.0834
$ws.Range("J113").Value = 1694.2
$ws.Range("K113").Value = 1604.0834
$ws.Range("L113").Value = 1694.2
$ws.Range("M113").Value = 565.9166
$ws.Range("N113").Value = -6034.2
# Row 122
$ws.Range("H122").Value = 4356.9375
$ws.Range("I122").Value = 4314.067
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 12942.201
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -10492.201
$ws.Range("N122").Value = -19900
# Row 123
$ws.Range("H123").Value = 21666.666
$ws.Range("J123").Value = 21666.666
$ws.Range("L123").Value = 21666.666
$ws.Range("N123").Value = -26566.666
# Row 126
$ws.Range("H126").Value = 2286.4285
$ws.Range("I126").Value = 1888.3334
$ws.Range("J126").Value = 2585
$ws.Range("K126").Value = 5665.0002
$ws.Range("L126").Value = 7755
$ws.Range("M126").Value = -3195.0002
$ws.Range("N126").Value = -12695
# Row 132
$ws.Range("H132").Value = 2813.4849
$ws.Range("I132").Value = 3172.1875
$ws.Range("J132").Value = 2475.8823
$ws.Range("K132").Value = 9516.5625
$ws.Range("L132").Value = 7427.646900000001
$ws.Range("M132").Value = -6986.5625
$ws.Range("N132").Value = -12487.6469

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2563.8
$ws.Range("I7").Value = 2374.2
$ws.Range("K7").Value = 2374.2
$ws.Range("M7").Value = -2262.2
# Row 22
$ws.Range("H22").Value = 1451
$ws.Range("J22").Value = 1451
$ws.Range("L22").Value = 1451
$ws.Range("N22").Value = -2041
# Row 27
$ws.Range("H27").Value = 1451
$ws.Range("J27").Value = 1451
$ws.Range("L27").Value = 1451
$ws.Range("N27").Value = -1665
# Row 55
$ws.Range("H55").Value = 209.9375
$ws.Range("I55").Value = 129.61539
$ws.Range("J55").Value = 558
$ws.Range("K55").Value = 129.61539
$ws.Range("L55").Value = 558
$ws.Range("M55").Value = 43.38461000000001
$ws.Range("N55").Value = -904
# Row 68
$ws.Range("H68").Value = 1208.1818
$ws.Range("I68").Value = 1237.5
$ws.Range("J68").Value = 1130
$ws.Range("K68").Value = 1237.5
$ws.Range("L68").Value = 1130
$ws.Range("M68").Value = -488.5
$ws.Range("N68").Value = -2628
# Row 71
$ws.Range("H71").Value = 1208.1818
$ws.Range("I71").Value = 1237.5
$ws.Range("J71").Value = 1130
$ws.Range("K71").Value = 6187.5
$ws.Range("L71").Value = 5650
$ws.Range("M71").Value = -2443.5
$ws.Range("N71").Value = -13138
# Row 94
$ws.Range("H94").Value = 4995
$ws.Range("J94").Value = 4995
$ws.Range("L94").Value = 4995
$ws.Range("N94").Value = -6347
# Row 126
$ws.Range("H126").Value = 2563.8
$ws.Range("I126").Value = 2374.2
$ws.Range("K126").Value = 7122.599999999999
$ws.Range("M126").Value = -4652.599999999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 15627236
$ws.Range("I122").Value = 19233190
$ws.Range("J122").Value = 1431.6666
$ws.Range("K122").Value = 57699570
$ws.Range("L122").Value = 4294.9998
$ws.Range("M122").Value = -57697120
$ws.Range("N122").Value = -9194.9998
